$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Insert a new "2022-Q3" sheet, positioned right before the existing
#    "2022-Q2" sheet (so the final tab order becomes:
#    总计, 2022-Q3, 2022-Q2, 2022-Q1).
#    We duplicate the "2022-Q2" sheet so the new sheet starts out with the
#    same layout/styling (headers, borders, column order) as the other
#    quarterly fund sheets, then we overwrite the data row with the new
#    Q3 figures.
# ---------------------------------------------------------------------------
$srcQ2 = $wb.Worksheets.Item("2022-Q2")
$srcQ2.Copy($srcQ2)
$newQ3 = $wb.Worksheets.Item("2022-Q2 (2)")
$newQ3.Name = "2022-Q3"

# Fund name changed slightly (QDII -> （QDII）) - plain text, no special
# handling needed, keeps its existing (default) formatting.
$newQ3.Range("C2").Value = "华夏全球科技先锋混合（QDII）"

# D2, E2, F2, G2 are stored as TEXT (not numbers) in the source data, e.g.
# "0.59" rather than 0.59. Excel auto-detects numeric-looking input and
# stores it as a number, so force the cell to Text format first and assign
# the value; then re-paste the (unstyled) formatting from B2 - a plain,
# default-style text cell in the same row - so no stray number-format /
# style index is left on the cell.
$newQ3.Range("B2").Copy()

$newQ3.Range("D2").NumberFormat = "@"
$newQ3.Range("D2").Value = "0.59"
$newQ3.Range("D2").PasteSpecial(-4122)

$newQ3.Range("E2").NumberFormat = "@"
$newQ3.Range("E2").Value = "86.79"
$newQ3.Range("E2").PasteSpecial(-4122)

$newQ3.Range("F2").NumberFormat = "@"
$newQ3.Range("F2").Value = "7.48"
$newQ3.Range("F2").PasteSpecial(-4122)

$newQ3.Range("G2").NumberFormat = "@"
$newQ3.Range("G2").Value = "0.0441"
$newQ3.Range("G2").PasteSpecial(-4122)

# H2 is a genuine number.
$newQ3.Range("H2").Value = 5

# ---------------------------------------------------------------------------
# 2) Update the "总计" (summary) sheet: the newest quarter becomes the
#    second row, the rest shift down, and a row is added for the oldest
#    quarter that is now falling out of the two-row window.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("B2").Value = "2022-Q3"
$total.Range("D2").Value = 0.04

$total.Range("B3").Value = "2022-Q2"
$total.Range("D3").Value = 0.03

# New row 4 (2022-Q1) - copy the style of row 3's "index" cell (A3, which
# carries the bordered/bold "index column" style) onto A4 before writing it.
$total.Range("A3").Copy()
$total.Range("A4").PasteSpecial(-4122)
$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 1
$total.Range("D4").Value = 0.03

# Restore the originally-selected/active tab ("2022-Q1" was the active sheet
# before the edit); the sheet Copy() operation above would otherwise leave
# the newly inserted "2022-Q3" sheet selected instead.
$wb.Worksheets.Item("2022-Q1").Activate()

